$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the Price (D) cells that are about to receive a new value as
# Text. Several new values look like plain decimals (e.g. "3.20", "1.00")
# and Excel would otherwise silently coerce them to numbers (dropping the
# trailing zero) the same way typing them into the grid would.
# (Applied as separate single-area Range calls -- multi-area "A,B" union
# strings only reliably format the first area in this host.)
$ws.Range("D2:D6").NumberFormat = "@"
$ws.Range("D9:D10").NumberFormat = "@"
$ws.Range("D12:D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32:D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38:D47").NumberFormat = "@"
$ws.Range("D49:D51").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = "51.125.57"
$ws.Cells.Item(2, 5).Value = "  +0.68%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.959.81"
$ws.Cells.Item(3, 5).Value = "  +1.31%  "

# Row 4
$ws.Cells.Item(4, 4).Value = "1.00"
$ws.Cells.Item(4, 5).Value = "  +0.12%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "379.81"
$ws.Cells.Item(5, 5).Value = "  +2.59%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "102.28"
$ws.Cells.Item(6, 5).Value = "  +2.42%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +2.69%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.04%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "0.588"
$ws.Cells.Item(9, 5).Value = "  +2.84%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "36.48"
$ws.Cells.Item(10, 5).Value = "  +2.69%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  -0.72%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "0.0856"
$ws.Cells.Item(12, 5).Value = "  +2.31%  "

# Row 13
$ws.Cells.Item(13, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(13, 4).Value = "3.424.10"
$ws.Cells.Item(13, 5).Value = "  +1.53%  "

# Row 14
$ws.Cells.Item(14, 2).Value = "Polkadot"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(14, 4).Value = "7.77"
$ws.Cells.Item(14, 5).Value = "  +6.05%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "18.25"
$ws.Cells.Item(15, 5).Value = "  +2.75%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "2.972.81"
$ws.Cells.Item(16, 5).Value = "  +2.25%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "11.11"
$ws.Cells.Item(17, 5).Value = "  +1.10%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "0.992"
$ws.Cells.Item(18, 5).Value = "  +4.43%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "51.177.63"
$ws.Cells.Item(19, 5).Value = "  +1.07%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "3.20"
$ws.Cells.Item(20, 5).Value = "  +3.14%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "12.50"
$ws.Cells.Item(21, 5).Value = "  +2.82%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "0.0₃0957"
$ws.Cells.Item(22, 5).Value = "  +1.37%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "70.08"
$ws.Cells.Item(23, 5).Value = "  +3.14%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "266.46"
$ws.Cells.Item(24, 5).Value = "  +1.79%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "3.19"
$ws.Cells.Item(25, 5).Value = "  +4.75%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "7.79"
$ws.Cells.Item(26, 5).Value = "  -3.28%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "7.42"
$ws.Cells.Item(27, 5).Value = "  +2.89%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  +0.00%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  +2.49%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "0.165"
$ws.Cells.Item(30, 5).Value = "  +1.72%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  +0.65%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "10.27"
$ws.Cells.Item(32, 5).Value = "  +5.06%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "34.56"
$ws.Cells.Item(33, 5).Value = "  +6.62%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "51.19"
$ws.Cells.Item(34, 5).Value = "  +1.66%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  +1.85%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "0.0436"
$ws.Cells.Item(36, 5).Value = "  +0.55%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  -0.12%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "3.25"
$ws.Cells.Item(38, 5).Value = "  +6.48%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "0.117"
$ws.Cells.Item(39, 5).Value = "  +2.50%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "Stacks"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(40, 4).Value = "2.54"
$ws.Cells.Item(40, 5).Value = "  +5.30%  "

# Row 41
$ws.Cells.Item(41, 2).Value = "ARBITRUM"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(41, 4).Value = "1.83"
$ws.Cells.Item(41, 5).Value = "  +4.55%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "Celestia"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Cells.Item(42, 4).Value = "16.51"
$ws.Cells.Item(42, 5).Value = "  +3.36%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "125.20"
$ws.Cells.Item(43, 5).Value = "  +5.30%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "3.56"
$ws.Cells.Item(44, 5).Value = "  +11.71%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "21.49"
$ws.Cells.Item(45, 5).Value = "  +4.63%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "WEMIXToken"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(46, 4).Value = "2.02"
$ws.Cells.Item(46, 5).Value = "  +0.04%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "ApeXProtocol"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Cells.Item(47, 4).Value = "2.37"
$ws.Cells.Item(47, 5).Value = "  +4.52%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  +0.12%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "2.028.80"
$ws.Cells.Item(49, 5).Value = "  +3.33%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "BEAM"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Cells.Item(50, 4).Value = "0.0321"
$ws.Cells.Item(50, 5).Value = "  +0.73%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "WOONetwork"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Cells.Item(51, 4).Value = "0.516"
$ws.Cells.Item(51, 5).Value = "  +14.27%  "
